# Update the date line at the top of the document.
$d = $word.ActiveDocument
$d.Content.Find.Execute("2026-02-20 Friday", $true, $false, $false, $false, $false, $true, 1, $false, "2026-02-21 Saturday", 2)

# Update the 25 division-fact cells (5 data rows x 5 columns) in the table.
# Row indices below are 1-based Word table rows; data rows are 1, 5, 9, 13, 17
# (rows 2-4, 6-8, 10-12, 14-16, 18-20 are blank spacer rows, left untouched).
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text  = "71÷7=10, 1"
$t.Cell(1,2).Range.Text  = "54÷8=6, 6"
$t.Cell(1,3).Range.Text  = "43÷5=8, 3"
$t.Cell(1,4).Range.Text  = "79÷3=26, 1"
$t.Cell(1,5).Range.Text  = "45÷3=15, 0"

$t.Cell(5,1).Range.Text  = "73÷8=9, 1"
$t.Cell(5,2).Range.Text  = "13÷9=1, 4"
$t.Cell(5,3).Range.Text  = "22÷6=3, 4"
$t.Cell(5,4).Range.Text  = "47÷5=9, 2"
$t.Cell(5,5).Range.Text  = "72÷8=9, 0"

$t.Cell(9,1).Range.Text  = "87÷4=21, 3"
$t.Cell(9,2).Range.Text  = "65÷2=32, 1"
$t.Cell(9,3).Range.Text  = "95÷7=13, 4"
$t.Cell(9,4).Range.Text  = "18÷7=2, 4"
$t.Cell(9,5).Range.Text  = "97÷2=48, 1"

$t.Cell(13,1).Range.Text = "64÷6=10, 4"
$t.Cell(13,2).Range.Text = "18÷5=3, 3"
$t.Cell(13,3).Range.Text = "18÷9=2, 0"
$t.Cell(13,4).Range.Text = "79÷4=19, 3"
$t.Cell(13,5).Range.Text = "42÷3=14, 0"

$t.Cell(17,1).Range.Text = "12÷8=1, 4"
$t.Cell(17,2).Range.Text = "83÷8=10, 3"
$t.Cell(17,3).Range.Text = "72÷9=8, 0"
$t.Cell(17,4).Range.Text = "32÷8=4, 0"
$t.Cell(17,5).Range.Text = "19÷2=9, 1"
